$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("A2").Value = "21331a0562"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = "Q1"
$ws.Range("D2").Value = "Q7"

# Add new row 3 values
$ws.Range("A3").Value = "21331a0569"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = "Q4"
$ws.Range("D3").Value = "Q10"
